$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Objetivos:) -------------------------------------------------
# Replace the long objectives paragraph with the teacher's name/ID.
$ws.Range("B10").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C10").Value = "8767640 - Eduardo Ferro dos Santos"

# --- Row 13 ---------------------------------------------------------------
# Previously unlabeled row holding "8767640 - Eduardo Ferro dos Santos" in B/C.
# Now becomes the "Programa resumido:" row, with "Semestral" as its value.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = $ws.Rows(10).RowHeight

# --- Row 14 -----------------------------------------------------------
# Was "Programa resumido:" with the short syllabus list in B/C.
# Becomes the "Short syllabus:" label only (no B/C content).
# Use Clear() (not just ClearContents) so the B14/C14 cells are removed
# entirely, matching rows that only ever had an A-column label.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# --- Row 15 -----------------------------------------------------------
# Was "Short syllabus:" only. Becomes "Programa:" with "01/01/2015" in B/C.
# Copy B8/C8 (already "01/01/2015", stored as text) instead of typing the
# literal, so Excel doesn't reinterpret it as a date value.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Rows(15).RowHeight = $ws.Rows(16).RowHeight

# --- Row 16 -----------------------------------------------------------
# Was "Programa:" with the long syllabus text in B/C.
# Becomes the "Syllabus:" label only (no B/C content).
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# --- Row 17 -----------------------------------------------------------
# Was "Syllabus:" only (120pt row). Becomes "Avaliação:" at default height.
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows(17).AutoFit()

# --- Row 18 -----------------------------------------------------------
# Was "Avaliação:" only. Becomes "Método:" with the teacher name in B/C
# (same text + style as the updated B10/C10, so copy from there).
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))
$ws.Rows(18).RowHeight = $ws.Rows(19).RowHeight

# --- Row 19 -----------------------------------------------------------
# Was "Método:" with the "Duas Notas..." grading text in B/C (unchanged).
# Becomes "Critério:".
$ws.Range("A19").Value = "Critério:"

# --- Row 20 -----------------------------------------------------------
# Was "Critério:" with "MF = (N1+ N2)/2" in B/C (unchanged).
# Becomes "Norma de recuperação:".
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21 -----------------------------------------------------------
# Was "Norma de recuperação:" with the "NF = ..." text in B/C (unchanged).
# Becomes "Bibliografia:". Row height grows from 60 to 120 (matches old row 22).
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# --- Row 22 -----------------------------------------------------------
# Was "Bibliografia:" with the long bibliography text in B/C.
# This row is dropped entirely.
$ws.Rows(22).Delete()
